# Lecture partielle de l'EDT M1 MIAGE.
# Shift each week-start date forward by 1096 days (3 years) and relabel the
# weekday text accordingly: every "lundi" (Monday) row becomes "vendredi"
# (Friday), and the one "vendredi" (Friday) row becomes "mardi" (Tuesday).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Liste")

# Rows in column A holding the week-start date, with column B carrying the
# corresponding weekday label.
$rows = @(2, 5, 8, 11, 14, 17, 20, 23, 26, 29, 32, 35)

foreach ($r in $rows) {
    $dateCell = $ws.Cells.Item($r, 1)
    $oldDate = $dateCell.Value2
    $dateCell.Value2 = $oldDate + 1096

    $labelCell = $ws.Cells.Item($r, 2)
    $label = $labelCell.Value2
    if ($label -eq "lundi") {
        $labelCell.Value2 = "vendredi"
    } elseif ($label -eq "vendredi") {
        $labelCell.Value2 = "mardi"
    }
}
